# Insert two new weekly data rows for "Camote" (Zapallo) at row 658,
# pushing the existing rows 658:705 down to 660:707, then fill in the
# values for the newly inserted rows 658 and 659.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 658 (shifts 658:705 -> 660:707)
$ws.Rows("658:659").Insert()

# --- Row 658: "1a (guarda)" ---
$ws.Range("A658").Value = 8
$ws.Range("B658").Value = "Terminal La Palmera de La Serena"
$ws.Range("C658").Value = "Coquimbo"
$ws.Range("D658").Value = 44746
$ws.Range("E658").Value = 4
$ws.Range("F658").Value = 100112045
$ws.Range("G658").Value = "Zapallo"
$ws.Range("H658").Value = "Camote"
$ws.Range("I658").Value = "1a (guarda)"
$ws.Range("J658").Value = 1720
$ws.Range("K658").Value = 750
$ws.Range("L658").Value = 800
$ws.Range("M658").Value = 775
$ws.Range("N658").Value = "$/kilo (volumen en unidades)"
$ws.Range("O658").Value = "Región de O'Higgins"
$ws.Range("P658").Value = 775
$ws.Range("Q658").Value = 1
$ws.Range("R658").Value = "Hortaliza"

# --- Row 659: "2a (guarda)" ---
$ws.Range("A659").Value = 8
$ws.Range("B659").Value = "Terminal La Palmera de La Serena"
$ws.Range("C659").Value = "Coquimbo"
$ws.Range("D659").Value = 44746
$ws.Range("E659").Value = 4
$ws.Range("F659").Value = 100112045
$ws.Range("G659").Value = "Zapallo"
$ws.Range("H659").Value = "Camote"
$ws.Range("I659").Value = "2a (guarda)"
$ws.Range("J659").Value = 920
$ws.Range("K659").Value = 650
$ws.Range("L659").Value = 700
$ws.Range("M659").Value = 675
$ws.Range("N659").Value = "$/kilo (volumen en unidades)"
$ws.Range("O659").Value = "Región de O'Higgins"
$ws.Range("P659").Value = 675
$ws.Range("Q659").Value = 1
$ws.Range("R659").Value = "Hortaliza"
